$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels (A1, B1) ---
$ws.Range("A1").Value = "readout noise [électrons]"
$ws.Range("B1").Value = "Frequency [pixel count]"

# --- Column widths ---
# (the runtime's ColumnWidth setter quantizes the stored XML width to
# 1/6-character increments, so these inputs are chosen as the closest
# values that land on the target stored widths of 25.83203125,
# 22.6640625, 15.5, 17.6640625, 17.6640625)
$ws.Columns.Item(1).ColumnWidth = 25
$ws.Columns.Item(2).ColumnWidth = 21.833333333333332
$ws.Columns.Item(4).ColumnWidth = 14.666666666666666
$ws.Columns.Item(5).ColumnWidth = 16.833333333333332
$ws.Columns.Item(6).ColumnWidth = 16.833333333333332

# --- Readout info block (D1:E4) ---
$ws.Range("D1").Value = "Readout info : "

$ws.Range("D2").Value = "mean : "
$ws.Range("E2").Formula = "=AVERAGE(A2:A43)"

$ws.Range("D3").Value = "stdev : "
$ws.Range("E3").Formula = "=STDEV(A2:A43)"

$ws.Range("D4").Value = "variance : "
$ws.Range("E4").Formula = "=E3^2"

# --- Selection / view state ---
[void]$ws.Range("L12").Select()
